$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20240805-102356-"

# Update the reference date (column G) for every data row from 45506 to 45509
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45509
}

# Update the handful of rows whose Saldo Previsto / Vl. Total values were corrected
$ws.Cells.Item(15, 5).Value = 999.9
$ws.Cells.Item(15, 8).Value = 999.9
$ws.Cells.Item(99, 5).Value = 982.54
$ws.Cells.Item(99, 8).Value = 982.54
$ws.Cells.Item(102, 5).Value = 726.16
$ws.Cells.Item(102, 8).Value = 726.16
$ws.Cells.Item(104, 5).Value = 515
$ws.Cells.Item(104, 8).Value = 515
$ws.Cells.Item(108, 5).Value = 985.55
$ws.Cells.Item(108, 8).Value = 985.55
$ws.Cells.Item(112, 5).Value = 34.15
$ws.Cells.Item(112, 8).Value = 34.15
$ws.Cells.Item(113, 5).Value = 14.07
$ws.Cells.Item(113, 8).Value = 14.07
$ws.Cells.Item(132, 5).Value = 985.04
$ws.Cells.Item(132, 8).Value = 985.04
$ws.Cells.Item(143, 5).Value = 999
$ws.Cells.Item(143, 8).Value = 999
$ws.Cells.Item(173, 5).Value = 973.16
$ws.Cells.Item(173, 8).Value = 973.16
$ws.Cells.Item(249, 5).Value = 973.46
$ws.Cells.Item(249, 8).Value = 973.46
$ws.Cells.Item(264, 5).Value = 988.82
$ws.Cells.Item(264, 8).Value = 988.82
$ws.Cells.Item(265, 5).Value = 980.05
$ws.Cells.Item(265, 8).Value = 980.05
$ws.Cells.Item(273, 5).Value = 984.88
$ws.Cells.Item(273, 8).Value = 984.88

# Clear the lingering cell selection from the previous save
$ws.Range("A1").Select()
